$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 2 values
$ws.Range("B2").Value = 81.11979141831398
$ws.Range("C2").Value = 81.11979141831398
$ws.Range("D2").Value = 74.58333298563957
$ws.Range("E2").Value = 66.71874970197678

# Add new header columns
$ws.Range("F1").Value = "epoch250"
$ws.Range("G1").Value = "epoch300"
$ws.Range("H1").Value = "epoch350"
$ws.Range("I1").Value = "epoch400"

# Add new row 2 values
$ws.Range("F2").Value = 63.59374970197678
$ws.Range("G2").Value = 61.24999970197678
$ws.Range("H2").Value = 58.51562470197678
$ws.Range("I2").Value = 58.51562470197678
